$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-31 Tuesday" "2023-11-01 Wednesday"

Replace-Text "85×46=3910" "36×50=1800"
Replace-Text "89×70=6230" "57×50=2850"
Replace-Text "51×21=1071" "21×27=567"
Replace-Text "37×24=888" "94×43=4042"
Replace-Text "50×33=1650" "45×65=2925"

Replace-Text "32×28=896" "52×50=2600"
Replace-Text "92×90=8280" "77×76=5852"
Replace-Text "68×20=1360" "65×29=1885"
Replace-Text "50×83=4150" "60×34=2040"
Replace-Text "29×27=783" "39×62=2418"

Replace-Text "39×72=2808" "69×12=828"
Replace-Text "25×77=1925" "65×94=6110"
Replace-Text "31×11=341" "31×41=1271"
Replace-Text "89×71=6319" "48×48=2304"
Replace-Text "14×27=378" "42×58=2436"

Replace-Text "46×88=4048" "31×17=527"
Replace-Text "11×76=836" "57×30=1710"
Replace-Text "22×11=242" "36×39=1404"
Replace-Text "20×92=1840" "92×39=3588"
Replace-Text "57×36=2052" "72×28=2016"

Replace-Text "82×64=5248" "68×16=1088"
Replace-Text "85×11=935" "53×81=4293"
Replace-Text "59×59=3481" "90×44=3960"
Replace-Text "43×44=1892" "16×91=1456"
Replace-Text "93×36=3348" "92×43=3956"

Write-Output "Done"
